# add save column in s_vals sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: copy the existing header formatting (bold font,
# border, centered alignment - same style used by B1:G1) from G1, then
# set its text to "Save".
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell H2 holding the numeric "Save" value for the row.
$ws.Range("H2").Value = 0
